# feat(excel2json)!: remove gui-attributes that do not have an effect in the APP (DEV-5482)
#
# The "gui_attributes" column (P) held several values that the app never
# actually reads, so they are cleared out. One gui_element value ("Radio")
# is also corrected to "List" (it is paired with a ListValue object, like
# every other "List" row, and "Radio" is not a real gui_element option).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Drop the now-meaningless gui_attributes entries.
$obsoleteGuiAttributes = @(
    "P2",   # hasAuthor            numprops: 1
    "P3",   # hasBibliographicReference  maxlength: 128, size: 128
    "P7",   # hasCopyright         maxlength: 128, size: 64
    "P8",   # hasCreator           numprops: 1
    "P14",  # hasImage             numprops: 1
    "P15",  # hasInventoryNumber   maxlength: 80, size: 25
    "P17",  # hasName              maxlength: 128, size: 128
    "P19",  # hasRelatedArtwork    numprops: 1
    "P25",  # hasWeight            maxlength: 255, size: 80
    "P26",  # inInstitution        numprops: 1
    "P28",  # hasIdentifier        maxlength: 128, size: 128
    "P29"   # hasChildren          max: 25.0, min: 0
)
foreach ($cellRef in $obsoleteGuiAttributes) {
    $ws.Range($cellRef).ClearContents()
}

# hasFlatList used the non-functional "Radio" gui_element; fix it to "List"
# to match how every other ListValue row is configured.
$ws.Range("O13").Value = "List"

# Move the viewport/selection back to the top-left of the sheet.
$null = $ws.Range("O36").Select()
